$d = $word.ActiveDocument

# --- Change 1: MODIS_AOD download time 35 17 -> 40 17 ---
$d.Content.Find.Execute(
    "35 17 * * * /home/fkaragulian/MODIS_AOD/MODIS_AOD_download_hdf_new.sh 1>/home/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "40 17 * * * /home/fkaragulian/MODIS_AOD/MODIS_AOD_download_hdf_new.sh 1>/home/", 2) | Out-Null

# --- Change 2: ECMWF 12gmt 0 7 -> 00 7 ---
$d.Content.Find.Execute(
    "0 7 * * * /home/fkaragulian/ECMWF_forecasts/ECMWF_forecasts_12gmt.sh 1>/home/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "00 7 * * * /home/fkaragulian/ECMWF_forecasts/ECMWF_forecasts_12gmt.sh 1>/home/", 2) | Out-Null

# --- Change 3: ECMWF 00gmt 0 17 -> 35 17 ---
$d.Content.Find.Execute(
    "0 17 * * * /home/fkaragulian/ECMWF_forecasts/ECMWF_forecasts_00gmt.sh 1>/home/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "35 17 * * * /home/fkaragulian/ECMWF_forecasts/ECMWF_forecasts_00gmt.sh 1>/home/", 2) | Out-Null

# --- Change 4: remove one of the two blank paragraphs before "30 8 ..." and
#     rewrite that line (comment it out, bump time, lowercase the script name) ---
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "30 8 * * * /home/fkaragulian/WRF_UAE/scripts/execute_Model.sh*") {
        $targetPara = $p
        break
    }
}
$prevPara = $targetPara.Previous()
$prevPrevPara = $prevPara.Previous()
# both $prevPara and $prevPrevPara should be the two consecutive blank paragraphs;
# delete just one of them (the one immediately before the target stays as separator
# after deletion is applied to the other one)
$prevPrevPara.Range.Delete() | Out-Null

$d.Content.Find.Execute(
    "30 8 * * * /home/fkaragulian/WRF_UAE/scripts/execute_Model.sh 1>/home/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "# 40 12 * * * /home/fkaragulian/WRF_UAE/scripts/execute_model.sh 1>/home/", 2) | Out-Null

# --- Change 5: insert a new crontab entry (submit_WRF_Chem.sh) right after the
#     paragraph that now reads "# 40 12 ... execute_model.sh ... chem-error.log",
#     before the following blank separator paragraph ---
$chemPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*execute_model.sh 1>/home/fkaragulian/log/chemrun.log*") {
        $chemPara = $p
        break
    }
}
$blankAfter = $chemPara.Next()
$insPt = $d.Range($blankAfter.Range.Start, $blankAfter.Range.Start)
$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>06 17 * * * /home/fkaragulian/WRF_UAE/scripts/submit_WRF_Chem.sh 1&gt;/home/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fkaragulian</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/log/chem_submit.log</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPt.InsertXML($newParaXml)

# --- Change 6: short_execute_Model.sh schedule "# 28 16 13 10 *" -> "# 41 21 * * *" ---
$d.Content.Find.Execute(
    "# 28 16 13 10 * /home/fkaragulian/WRF_UAE/scripts/short_execute_Model.sh 1>/home/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "# 41 21 * * * /home/fkaragulian/WRF_UAE/scripts/short_execute_Model.sh 1>/home/", 2) | Out-Null

# --- Change 7: merge the trailing bookmark-only paragraph into the previous
#     paragraph (remove the paragraph break between them) ---
$n = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs($n - 1)
$markRange = $d.Range($secondToLast.Range.End - 1, $secondToLast.Range.End)
$markRange.Delete() | Out-Null
